# "Almost doen with Today Forecast"
# Adds the weather-icon filename lookup (column E) to the "images & icon"
# sheet, one icon per OpenWeatherMap condition-code group, and updates the
# saved selections on both sheets.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("images & icon")

# Map each data row on the "images & icon" sheet to its weather-icon file.
$iconMap = @{
    3  = "thunderstorm.png"
    4  = "thunderstorm.png"
    5  = "thunderstorm.png"
    6  = "thunderstorm.png"
    7  = "thunderstorm.png"
    8  = "thunderstorm.png"
    9  = "thunderstorm.png"
    10 = "thunderstorm.png"
    11 = "thunderstorm.png"
    12 = "thunderstorm.png"

    14 = "drizzle.png"
    15 = "drizzle.png"
    16 = "drizzle.png"
    17 = "drizzle.png"
    18 = "drizzle.png"
    19 = "drizzle.png"
    20 = "drizzle.png"
    21 = "drizzle.png"
    22 = "drizzle.png"

    24 = "rain.png"
    25 = "rain.png"
    26 = "rain.png"
    27 = "rain.png"
    28 = "rain.png"
    29 = "freezing-rain.png"
    30 = "rain.png"
    31 = "rain.png"
    32 = "rain.png"
    33 = "rain.png"

    35 = "snow.png"
    36 = "snow.png"
    37 = "snow.png"
    38 = "snow.png"
    39 = "snow.png"
    40 = "snow.png"
    41 = "snow.png"
    42 = "snow.png"
    43 = "snow.png"
    44 = "snow.png"
    45 = "snow.png"

    47 = "fog.png"
    48 = "fog.png"
    49 = "fog.png"
    50 = "fog.png"
    51 = "fog.png"
    52 = "fog.png"
    53 = "fog.png"
    54 = "fog.png"
    55 = "fog.png"
    56 = "fog.png"

    58 = "sun.png"

    60 = "clouds.png"
    61 = "clouds.png"
    62 = "clouds.png"
    63 = "clouds.png"
}

foreach ($row in ($iconMap.Keys | Sort-Object)) {
    $ws2.Range("E$row").Value = $iconMap[$row]
}

# Let the new column size itself to its contents, like Excel's own
# auto-fit would do after typing the values in.
$ws2.Columns("E:E").AutoFit() | Out-Null

# Restore the selections that were active when the workbook was last saved.
$ws1.Activate()
$ws1.Range("A20").Select() | Out-Null

$ws2.Activate()
$ws2.Range("A47:A56").Select() | Out-Null
